# Add two new columns, I (I0) and J (IF), to the sheet.
# I column is a constant of 1 for every data row except row 30 (which is 7).
# J column mirrors the existing H (IP) column value, except row 30 (which is 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the existing header cell H1 onto the two new header
# cells so they pick up the same bold/centered/bordered style, then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (rows 2-31) ---
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 30) {
        $iVal = 7
        $jVal = 8
    } else {
        $iVal = 1
        $jVal = $hVal
    }

    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
